# Apply naming-convention changes described in the commit message:
#  - T1..T5 header cells become T_1..T_5
#  - NT1..NT5 header cells become NT_1..NT_5
#  - "CpG_Array" label (A22) becomes "Mean_beta-value"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header row (row 1): T1..T5 -> T_1..T_5, NT1..NT5 -> NT_1..NT_5
$ws.Range("B1").Value = "T_1"
$ws.Range("C1").Value = "T_2"
$ws.Range("D1").Value = "T_3"
$ws.Range("E1").Value = "T_4"
$ws.Range("F1").Value = "T_5"
$ws.Range("G1").Value = "NT_1"
$ws.Range("H1").Value = "NT_2"
$ws.Range("I1").Value = "NT_3"
$ws.Range("J1").Value = "NT_4"
$ws.Range("K1").Value = "NT_5"

# Row 22 label rename
$ws.Range("A22").Value = "Mean_beta-value"

# Update the active selection to mirror the saved workbook state
$ws.Range("A23").Select()
